$wb = $excel.ActiveWorkbook

$about = $wb.Worksheets.Item("About")
$data = $wb.Worksheets.Item("FoOMCtiL")

# --- "About" sheet: insert Notes section above Procedure, shifting rows 9-11 down to 12-14 ---
$about.Rows.Item(9).Resize(3, 1).Insert() | Out-Null

$about.Range("A9").Value = "Notes:"
$about.Range("A10").Value = "This variable sets the fraction of O&M costs in the electricity sector that is labor."

$about.Range("A9").Select() | Out-Null

# --- "FoOMCtiL" sheet: update header label & wrap text, reorder shared strings naturally follows ---
$data.Activate() | Out-Null
$data.Range("B1").Value = "Frac of O&M Costs (dimensionless)"
$data.Range("B1").WrapText = $true
$data.Rows.Item(1).RowHeight = 30
$data.Range("B1").Select() | Out-Null

# Leave "About" as the active / tab-selected sheet, with A11 as its last selection
$about.Activate() | Out-Null
$about.Range("A11").Select() | Out-Null
